# updated reliability & efficiency metrics
$wb = $excel.ActiveWorkbook

$wsEff = $wb.Worksheets.Item("Efficiency")
$wsRel = $wb.Worksheets.Item("Reliability")

# --- Efficiency sheet: refresh the latest (3rd) measurement block (rows 20-21) ---
$wsEff.Range("D20").Value = "42 / 2.7 KB"
$wsEff.Range("E20").Value = "49 / 1.8 KB"
$wsEff.Range("F20").Value = "48 / 3.0 KB"
$wsEff.Range("G20").Value = "49 / 1.8 KB"

$wsEff.Range("D21").Value = "1 / 64.3 B"
$wsEff.Range("E21").Value = "1 / 50 B"
$wsEff.Range("F21").Value = "1 / 62.5 B"
$wsEff.Range("G21").Value = "1 / 50 B"

# --- Reliability sheet: new 4th date block (col S) + new Issue-close-time series (cols Q:T) ---
$wsRel.Range("S3").Value = 42709

$wsRel.Range("Q8").Value = "#446"
$wsRel.Range("R8").Value = 33

$wsRel.Range("Q9").Value = "#445"
$wsRel.Range("R9").Value = 30

$wsRel.Range("Q10").Value = "#443"
$wsRel.Range("R10").Value = 226
$wsRel.Range("S10").ClearContents()

$wsRel.Range("Q11").Value = "#429"
$wsRel.Range("R11").Value = 202

$wsRel.Range("Q12").Value = "#424"
$wsRel.Range("R12").Value = 164

$wsRel.Range("Q13").Value = "#422"
$wsRel.Range("R13").Value = 163

$wsRel.Range("Q14").Value = "#420"
$wsRel.Range("R14").Value = 156

$wsRel.Range("Q15").Value = "#411"
$wsRel.Range("R15").Value = 178

$wsRel.Range("Q16").Value = "#400"
$wsRel.Range("R16").Value = 166

$wsRel.Range("Q17").Value = "#386"
$wsRel.Range("R17").Value = 161

$wsRel.Range("Q18").Value = "#385"
$wsRel.Range("R18").Value = 191

$wsRel.Range("Q19").Value = "#384"
$wsRel.Range("R19").Value = 182

$wsRel.Range("Q20").Value = "#382"
$wsRel.Range("R20").Value = 161

$wsRel.Range("Q21").Value = "#381"
$wsRel.Range("R21").Value = 178

$wsRel.Range("Q22").Value = "#380"
$wsRel.Range("R22").Value = 168

$wsRel.Range("Q23").Value = "#379"
$wsRel.Range("R23").Value = 163

$wsRel.Range("Q24").Value = "#378"
$wsRel.Range("R24").Value = 182

$wsRel.Range("Q25").Value = "#377"
$wsRel.Range("R25").Value = 158

$wsRel.Range("Q26").Value = "#375"
$wsRel.Range("R26").Value = 181

$wsRel.Range("Q27").Value = "#370"
$wsRel.Range("R27").Value = 237

$wsRel.Range("Q28").Value = "#359"
$wsRel.Range("R28").Value = 199

$wsRel.Range("T8").Formula = "=AVERAGE(R8:R28)"

# New column (T) needs the same kind of width as the other data columns
$wsRel.Columns.Item(20).ColumnWidth = 27

# --- View / selection state ---
# Efficiency tab becomes inactive; its selection moves to F21
$wsEff.Range("F21").Select()

# Reliability tab becomes the active tab; selection moves to the new T8 cell
$wsRel.Range("T8").Select()
